# Scheduled runner update: refresh computed market-price / profit columns (H-N)
# across the per-job Leve profitability sheets.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1835.75
$ws.Range("I111").Value = 2297
$ws.Range("J111").Value = 1242.7142
$ws.Range("K111").Value = 6891
$ws.Range("L111").Value = 3728.1426
$ws.Range("M111").Value = -3824
$ws.Range("N111").Value = -9862.142599999999
$ws.Range("H116").Value = 7000
$ws.Range("I116").Value = 7000
$ws.Range("K116").Value = 7000
$ws.Range("M116").Value = -3558
$ws.Range("H125").Value = 27781436
$ws.Range("I125").Value = 41670230
$ws.Range("J125").Value = 3845.3333
$ws.Range("K125").Value = 375032070
$ws.Range("L125").Value = 34607.9997
$ws.Range("M125").Value = -375029610
$ws.Range("N125").Value = -39527.9997
$ws.Range("H132").Value = 2646.2666
$ws.Range("I132").Value = 3474.375
$ws.Range("J132").Value = 1699.8572
$ws.Range("K132").Value = 10423.125
$ws.Range("L132").Value = 5099.571599999999
$ws.Range("M132").Value = -7893.125
$ws.Range("N132").Value = -10159.5716
$ws.Range("H138").Value = 11192.75
$ws.Range("J138").Value = 11374.675
$ws.Range("L138").Value = 34124.02499999999
$ws.Range("N138").Value = -44404.02499999999

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1053.909
$ws.Range("I2").Value = 1029.3
$ws.Range("K2").Value = 1029.3
$ws.Range("M2").Value = -916.3
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H116").Value = 1053.909
$ws.Range("I116").Value = 1029.3
$ws.Range("K116").Value = 1029.3
$ws.Range("M116").Value = 1264.7

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1053.909
$ws.Range("I3").Value = 1029.3
$ws.Range("K3").Value = 1029.3
$ws.Range("M3").Value = -915.3
$ws.Range("H20").Value = 3523.3845
$ws.Range("I20").Value = 2178.6
$ws.Range("K20").Value = 2178.6
$ws.Range("M20").Value = -1931.6
$ws.Range("H99").Value = 1057.8334
$ws.Range("I99").Value = 199
$ws.Range("J99").Value = 1916.6666
$ws.Range("K99").Value = 199
$ws.Range("L99").Value = 1916.6666
$ws.Range("M99").Value = 1299
$ws.Range("N99").Value = -4912.6666

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2906
$ws.Range("J10").Value = 5002.6665
$ws.Range("L10").Value = 5002.6665
$ws.Range("N10").Value = -5280.6665
$ws.Range("H16").Value = 2020.625
$ws.Range("I16").Value = 637.6
$ws.Range("K16").Value = 637.6
$ws.Range("M16").Value = -350.6
$ws.Range("H32").Value = 5750
$ws.Range("I32").Value = 5750
$ws.Range("K32").Value = 5750
$ws.Range("M32").Value = -5434
$ws.Range("H113").Value = 2020.625
$ws.Range("I113").Value = 637.6
$ws.Range("K113").Value = 637.6
$ws.Range("M113").Value = 1532.4
$ws.Range("H132").Value = 9483.777
$ws.Range("I132").Value = 3967.5557
$ws.Range("K132").Value = 11902.6671
$ws.Range("M132").Value = -9372.667099999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 140.90909
$ws.Range("I4").Value = 140.90909
$ws.Range("K4").Value = 422.72727
$ws.Range("M4").Value = -310.72727

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2005
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3849
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 33335502
$ws.Range("I70").Value = 33335502
$ws.Range("K70").Value = 33335502
$ws.Range("M70").Value = -33335232
$ws.Range("H73").Value = 33335502
$ws.Range("I73").Value = 33335502
$ws.Range("K73").Value = 33335502
$ws.Range("M73").Value = -33334566
$ws.Range("H80").Value = 3883
$ws.Range("I80").Value = 3841.75
$ws.Range("K80").Value = 3841.75
$ws.Range("M80").Value = -2843.75
$ws.Range("H83").Value = 3883
$ws.Range("I83").Value = 3841.75
$ws.Range("K83").Value = 19208.75
$ws.Range("M83").Value = -14216.75
$ws.Range("H97").Value = 1377
$ws.Range("I97").Value = 1003.6
$ws.Range("J97").Value = 1999.3334
$ws.Range("K97").Value = 1003.6
$ws.Range("L97").Value = 1999.3334
$ws.Range("M97").Value = -507.6
$ws.Range("N97").Value = -2991.3334
$ws.Range("H102").Value = 1063.8
$ws.Range("I102").Value = 1063.8
$ws.Range("K102").Value = 1063.8
$ws.Range("M102").Value = 558.2
$ws.Range("H132").Value = 2623.3157
$ws.Range("I132").Value = 2149.5881
$ws.Range("K132").Value = 6448.7643
$ws.Range("M132").Value = -3918.7643

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2634.5
$ws.Range("I2").Value = 845.4545000000001
$ws.Range("J2").Value = 6570.4
$ws.Range("K2").Value = 845.4545000000001
$ws.Range("L2").Value = 6570.4
$ws.Range("M2").Value = -733.4545000000001
$ws.Range("N2").Value = -6794.4
$ws.Range("H19").Value = 4337.125
$ws.Range("I19").Value = 582
$ws.Range("K19").Value = 582
$ws.Range("M19").Value = -412
$ws.Range("H32").Value = 2000
$ws.Range("J32").Value = 2000
$ws.Range("L32").Value = 2000
$ws.Range("N32").Value = -2634
$ws.Range("H74").Value = 69999
$ws.Range("J74").Value = 69999
$ws.Range("L74").Value = 69999
$ws.Range("N74").Value = -71995
$ws.Range("H77").Value = 69999
$ws.Range("J77").Value = 69999
$ws.Range("L77").Value = 209997
$ws.Range("N77").Value = -219981
$ws.Range("H136").Value = 4555.5557
$ws.Range("J136").Value = 4800
$ws.Range("L136").Value = 14400
$ws.Range("N136").Value = -19500

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 780.5
$ws.Range("J2").Value = 1550
$ws.Range("L2").Value = 1550
$ws.Range("N2").Value = -1774
$ws.Range("H75").Value = 24166.334
$ws.Range("J75").Value = 23999.8
$ws.Range("L75").Value = 23999.8
$ws.Range("N75").Value = -25871.8
$ws.Range("H78").Value = 24166.334
$ws.Range("J78").Value = 23999.8
$ws.Range("L78").Value = 71999.39999999999
$ws.Range("N78").Value = -81359.39999999999
$ws.Range("H87").Value = 58940
$ws.Range("J87").Value = 58940
$ws.Range("L87").Value = 58940
$ws.Range("N87").Value = -61436
$ws.Range("H90").Value = 58940
$ws.Range("J90").Value = 58940
$ws.Range("L90").Value = 176820
$ws.Range("N90").Value = -189300
$ws.Range("H132").Value = 3759
$ws.Range("I132").Value = 3283.6
$ws.Range("K132").Value = 9850.799999999999
$ws.Range("M132").Value = -7320.799999999999

